$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7, 6, 5 (from bottom up) so remaining row indices don't shift
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# Row 2
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Bdnf"
$ws.Range("C2").Value = "Ngfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9759805
$ws.Range("H2").Value = 1.951961
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.337665
$ws.Range("N2").Value = 0.67533
$ws.Range("O2").Value = 0.01438349055657064
$ws.Range("P2").Value = 0.0143657952272707
$ws.Range("Q2").Value = 0.3295544555325
$ws.Range("R2").Value = 1.31821782213
$ws.Range("S2").Value = 0.01438349055657064
$ws.Range("T2").Value = 0.0143657952272707

# Row 3
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Bdnf"
$ws.Range("C3").Value = "Ngfr"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9759805
$ws.Range("H3").Value = 1.951961
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05783366666666667
$ws.Range("N3").Value = 0.173501
$ws.Range("O3").Value = 0.002463536340313192
$ws.Range("P3").Value = 0.003690758351808291
$ws.Range("Q3").Value = 0.05644453091016666
$ws.Range("R3").Value = 0.338667185461
$ws.Range("S3").Value = 0.002463536340313192
$ws.Range("T3").Value = 0.003690758351808291

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Bdnf"
$ws.Range("C4").Value = "Ngfr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9759805
$ws.Range("H4").Value = 1.951961
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 23.0803745
$ws.Range("N4").Value = 46.160749
$ws.Range("O4").Value = 0.9831529731031161
$ws.Range("P4").Value = 0.981943446420921
$ws.Range("Q4").Value = 22.52599544469725
$ws.Range("R4").Value = 90.10398177878899
$ws.Range("S4").Value = 0.9831529731031161
$ws.Range("T4").Value = 0.981943446420921
